# Add a new "path" row (row 6) to the Sprites table, mirroring the layout
# of the existing rows (column A = label, columns B:I = per-sprite values
# in ID order: p1, e1, e2, b1, b2, s1, s2, r1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "path"
$ws.Range("B6").Value = "sprites/player/p1"
$ws.Range("C6").Value = "sprites/enemy/e1"
$ws.Range("D6").Value = "sprites/enemy/e2"
$ws.Range("E6").Value = "blocks/b1"
$ws.Range("F6").Value = "blocks/b2"
$ws.Range("G6").Value = "blocks/s1"
$ws.Range("H6").Value = "blocks/s2"
$ws.Range("I6").Value = "blocks/r1"

$ws.Range("H6").Select()
